$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '67.635.25'
Set-TextValue $ws.Range("E2") '  -7.65%  '

Set-TextValue $ws.Range("D3") '3.683.26'
Set-TextValue $ws.Range("E3") '  -7.51%  '

Set-TextValue $ws.Range("D4") '0.995'
Set-TextValue $ws.Range("E4") '  -0.58%  '

Set-TextValue $ws.Range("D5") '567.50'
Set-TextValue $ws.Range("E5") '  -6.77%  '

Set-TextValue $ws.Range("D6") '172.55'
Set-TextValue $ws.Range("E6") '  +0.11%  '

Set-TextValue $ws.Range("D7") '3.661.53'
Set-TextValue $ws.Range("E7") '  -7.85%  '

Set-TextValue $ws.Range("D8") '0.621'
Set-TextValue $ws.Range("E8") '  -9.74%  '

Set-TextValue $ws.Range("D9") '1.00'
Set-TextValue $ws.Range("E9") '  +0.09%  '

Set-TextValue $ws.Range("D10") '0.700'
Set-TextValue $ws.Range("E10") '  -11.72%  '

Set-TextValue $ws.Range("D11") '0.161'
Set-TextValue $ws.Range("E11") '  -13.38%  '

Set-TextValue $ws.Range("D12") '50.95'
Set-TextValue $ws.Range("E12") '  -11.38%  '

Set-TextValue $ws.Range("D13") '0.0000291'
Set-TextValue $ws.Range("E13") '  -14.15%  '

Set-TextValue $ws.Range("D14") '10.41'
Set-TextValue $ws.Range("E14") '  -11.43%  '

Set-TextValue $ws.Range("D15") '4.225.14'
Set-TextValue $ws.Range("E15") '  -8.59%  '

Set-TextValue $ws.Range("D16") '3.652.92'

$ws.Range("B17").Value = 'Chainlink'
$ws.Range("C17").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D17") '19.25'
Set-TextValue $ws.Range("E17") '  -8.14%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range("D18") '0.126'
Set-TextValue $ws.Range("E18") '  -3.51%  '

Set-TextValue $ws.Range("D19") '12.80'
Set-TextValue $ws.Range("E19") '  -10.65%  '

Set-TextValue $ws.Range("D20") '1.12'
Set-TextValue $ws.Range("E20") '  -11.07%  '

Set-TextValue $ws.Range("D21") '66.951.07'
Set-TextValue $ws.Range("E21") '  -8.56%  '

Set-TextValue $ws.Range("D22") '403.70'
Set-TextValue $ws.Range("E22") '  -13.98%  '

Set-TextValue $ws.Range("D23") '4.42'
Set-TextValue $ws.Range("E23") '  -8.20%  '

Set-TextValue $ws.Range("D24") '87.26'
Set-TextValue $ws.Range("E24") '  -9.63%  '

Set-TextValue $ws.Range("D25") '3.02'
Set-TextValue $ws.Range("E25") '  -11.57%  '

Set-TextValue $ws.Range("D26") '12.63'
Set-TextValue $ws.Range("E26") '  -11.74%  '

Set-TextValue $ws.Range("D27") '10.60'
Set-TextValue $ws.Range("E27") '  -5.46%  '

Set-TextValue $ws.Range("E28") '  +0.33%  '

Set-TextValue $ws.Range("D29") '3.69'
Set-TextValue $ws.Range("E29") '  -13.48%  '

Set-TextValue $ws.Range("D30") '9.36'
Set-TextValue $ws.Range("E30") '  -12.34%  '

Set-TextValue $ws.Range("D31") '32.41'
Set-TextValue $ws.Range("E31") '  -11.13%  '

Set-TextValue $ws.Range("D32") '7.54'
Set-TextValue $ws.Range("E32") '  -6.92%  '

Set-TextValue $ws.Range("D33") '12.39'
Set-TextValue $ws.Range("E33") '  -11.89%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D34") '0.115'
Set-TextValue $ws.Range("E34") '  -11.57%  '

$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Range("D35") '64.61'
Set-TextValue $ws.Range("E35") '  -8.47%  '

Set-TextValue $ws.Range("D36") '42.84'
Set-TextValue $ws.Range("E36") '  -14.22%  '

$ws.Range("B37").Value = 'PEPE'
$ws.Range("C37").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Range("D37") '0.0₃0893'
Set-TextValue $ws.Range("E37") '  -13.15%  '

$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D38") '576.85'
Set-TextValue $ws.Range("E38") '  -10.00%  '

Set-TextValue $ws.Range("E39") '  -0.17%  '

$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws.Range("D40") '0.391'
Set-TextValue $ws.Range("E40") '  -9.73%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D41") '0.992'
Set-TextValue $ws.Range("E41") '  -0.98%  '

Set-TextValue $ws.Range("D42") '0.132'
Set-TextValue $ws.Range("E42") '  -11.04%  '

Set-TextValue $ws.Range("D43") '2.96'
Set-TextValue $ws.Range("E43") '  -8.75%  '

Set-TextValue $ws.Range("D44") '2.95'
Set-TextValue $ws.Range("E44") '  -13.63%  '

Set-TextValue $ws.Range("D45") '0.0433'
Set-TextValue $ws.Range("E45") '  -10.99%  '

Set-TextValue $ws.Range("D46") '2.54'
Set-TextValue $ws.Range("E46") '  -3.04%  '

Set-TextValue $ws.Range("D47") '9.07'
Set-TextValue $ws.Range("E47") '  -14.45%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D48") '0.133'
Set-TextValue $ws.Range("E48") '  -11.11%  '

$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D49") '3.15'
Set-TextValue $ws.Range("E49") '  -8.35%  '

Set-TextValue $ws.Range("D50") '2.66'
Set-TextValue $ws.Range("E50") '  -5.87%  '

Set-TextValue $ws.Range("D51") '2.691.35'
Set-TextValue $ws.Range("E51") '  -4.57%  '
